$d = $word.ActiveDocument

# Replace $oldText with $newText, but scope the Find/Replace operation to a
# narrow Range computed from the plain-text offsets so that we do not
# accidentally touch / reformat neighbouring runs that happen to carry
# different character formatting (e.g. the literal "<tag>" markers which
# are styled in Courier New / blue in this document).
#
# $searchFrom (optional) disambiguates which occurrence to use, by
# restricting the search to start at or after the first occurrence of
# $searchFrom in the document text.
function Replace-Scoped {
    param(
        [string]$oldText,
        [string]$newText,
        [string]$searchFrom = $null
    )

    $full = $d.Content.Text
    $startSearch = 0
    if ($searchFrom) {
        $startSearch = $full.IndexOf($searchFrom)
        if ($startSearch -lt 0) {
            throw "Anchor not found: $searchFrom"
        }
    }

    $idx = $full.IndexOf($oldText, $startSearch)
    if ($idx -lt 0) {
        throw "Could not find text: $oldText"
    }

    $start = $idx
    $end = $idx + $oldText.Length
    $r = $d.Range($start, $end)
    $r.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                     $true, 1, $false, $newText, 2) | Out-Null
}

# 1. "the false one" -> "piercing"  (inside "...battery <add>&amp; the false one</add> are...")
Replace-Scoped "the false one" "piercing"

# 2. "ball " -> "ball, "  (in "Their ball of the <ms>King's caliber</ms>...")
Replace-Scoped "ball " "ball, " "Their "

# 3. " and battery weighs 30 " -> " &amp; of battery, is 30 "
Replace-Scoped " and battery weighs 30 " " &amp; of battery, is 30 "

# 4. "...is therefore lighter...because fifteen " -> "...thus lighter...for fifteen "
Replace-Scoped `
    "s and is therefore lighter than the canon's one. Therefore, It does not carry so many munitions because fifteen " `
    "s and thus lighter than that of the cannon. And by thus, it does not carry so much munition for fifteen "

# 5. canon/culverine comparison sentence rewrite
Replace-Scoped `
    "s are enough for its load. The canon has a bigger mouth due to the size of its cannonball, but the colverine is more precise and goes faster, having greater range due to its length. Its bre" `
    "s are enough for its load. The cannon has a bigger mouth due to the size of its ball, but the culverine is more vigorous &amp; is faster, having greater power due to its length. At its bre"

# 6. breech / culverines usage sentence rewrite
Replace-Scoped `
    "ch is two bullets and a third thick, the front is a bullet thick. Culverines are used for fighting fortifications from far away when it is not possible to easily approach them. And then canons come closer. They are used also to support the battery. Fifteen or sixteen " `
    "ch it is two balls &amp; one third thick, the front is one ball &amp; two thirds thick. Culverines are used for battering defenses from afar when one cannot easily make an approach. And cannons can come closer. They are used also to support the battery. Fifteen or sixteen "

# 7. horses / alloy sentence rewrite (also drops the "<sup>of artillery</sup>" markup)
Replace-Scoped `
    " are necessary to carry it. They are made of the same metal alloy than the canon, like all others smaller pieces <sup>of artillery</sup>. For these, we add a little bit more metal in order to make the melting run better. And for two " `
    " are needed to bring it. They are of the same alloy as the cannon, as are all pieces smaller than the moyen, for to these, we add a little bit more metal in order to make the melting run better. And for two "

# 8. "bastarde culverine" -> "bastarde coleverine"
Replace-Scoped "bastarde culverine" "bastarde coleverine"
